$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" text cell (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.0 = 27801.12 pesos`n✅ 27801.12 pesos = 6.97 = 967.14 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the rate cells N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 142.8
$wsTasas.Range("O10").Value = 3970
$wsTasas.Range("N12").Value = 3989.9
$wsTasas.Range("O12").Value = 138.8
